# GSDA-310 data sheet rework:
# - new header row (NOME / DATA / CPF / ENDEREÇO / CTPS)
# - date values move from column C to column B
# - new CPF column (text-formatted) plus ENDEREÇO and CTPS columns
# - second data row changes from "ihsa" to "Maria"
# - column widths best-fit to their new contents
# - selection left on C7

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (headers) ---
$ws.Range("A1").Value = "NOME"
$ws.Range("B1").Value = "DATA"
$ws.Range("C1").NumberFormat = "@"
$ws.Range("C1").Formula = "CPF"
$ws.Range("D1").Value = "ENDEREÇO"
$ws.Range("E1").Value = "CTPS"

# --- Row 2 ---
$ws.Range("A2").Value = "lucas"
$ws.Range("B2").Value2 = 45575
$ws.Range("B2").NumberFormat = "d-mmm"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Formula = "09706809700"
$ws.Range("D2").Value = "rua Graccho rangel, 553"
$ws.Range("E2").Value = "465468 123/RJ"

# --- Row 3 ---
$ws.Range("A3").Value = "Maria"
$ws.Range("B3").Value2 = 45576
$ws.Range("B3").NumberFormat = "d-mmm"
$ws.Range("C3").Value2 = 12345678910
$ws.Range("C3").NumberFormat = "@"
$ws.Range("D3").Value = "rua expedito sauasuh"
$ws.Range("E3").Value = "123456 64/RJ"

# --- Column widths (best fit to new content) ---
$ws.Columns.Item(1).ColumnWidth = 5.667
$ws.Columns.Item(2).ColumnWidth = 6
$ws.Columns.Item(3).ColumnWidth = 11.667
$ws.Columns.Item(4).ColumnWidth = 20.833
$ws.Columns.Item(5).ColumnWidth = 12.667

# --- Selection ---
$ws.Range("C7").Select() | Out-Null
